# Update the "想去人数" (want-to-go count) column F across the four sheets
# to reflect newly generated output numbers.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 306
$ws1.Range("F6").Value = 395
$ws1.Range("F7").Value = 853
$ws1.Range("F8").Value = 42
$ws1.Range("F9").Value = 506
$ws1.Range("F12").Value = 132
$ws1.Range("F14").Value = 230
$ws1.Range("F15").Value = 32
$ws1.Range("F17").Value = 6596
$ws1.Range("F18").Value = 62
$ws1.Range("F21").Value = 7535
$ws1.Range("F24").Value = 3384
$ws1.Range("F25").Value = 25
$ws1.Range("F26").Value = 1477
$ws1.Range("F27").Value = 884
$ws1.Range("F28").Value = 4512
$ws1.Range("F29").Value = 28
$ws1.Range("F31").Value = 67
$ws1.Range("F32").Value = 208
$ws1.Range("F34").Value = 1620
$ws1.Range("F36").Value = 150
$ws1.Range("F37").Value = 53
$ws1.Range("F38").Value = 17
$ws1.Range("F39").Value = 1190
$ws1.Range("F40").Value = 1715
$ws1.Range("F41").Value = 2130
$ws1.Range("F42").Value = 9

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 49

# --- Sheet "本地生活" (sheet3) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 241
$ws3.Range("F4").Value = 74

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 241
$ws4.Range("F5").Value = 74
$ws4.Range("F7").Value = 306
$ws4.Range("F8").Value = 395
$ws4.Range("F9").Value = 853
$ws4.Range("F10").Value = 42
$ws4.Range("F11").Value = 506
$ws4.Range("F15").Value = 132
$ws4.Range("F18").Value = 230
$ws4.Range("F19").Value = 32
$ws4.Range("F21").Value = 6596
$ws4.Range("F22").Value = 62
$ws4.Range("F25").Value = 7535
$ws4.Range("F28").Value = 3384
$ws4.Range("F29").Value = 25
$ws4.Range("F30").Value = 1477
$ws4.Range("F31").Value = 884
$ws4.Range("F32").Value = 4512
$ws4.Range("F33").Value = 29
$ws4.Range("F35").Value = 67
$ws4.Range("F36").Value = 49
$ws4.Range("F37").Value = 208
$ws4.Range("F39").Value = 1620
$ws4.Range("F41").Value = 150
$ws4.Range("F42").Value = 53
$ws4.Range("F43").Value = 17
$ws4.Range("F44").Value = 1190
$ws4.Range("F45").Value = 1715
$ws4.Range("F47").Value = 2130
$ws4.Range("F48").Value = 9
